$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '91.812.32'
$ws.Range("E2").Value = '  +3.07%  '

# Row 3
$ws.Range("D3").Value = '3.246.14'
$ws.Range("E3").Value = '  -1.19%  '

# Row 4
$ws.Range("E4").Value = '  +0.23%  '

# Row 5
$ws.Range("D5").Value = '218.13'
$ws.Range("E5").Value = '  +2.18%  '

# Row 6
$ws.Range("D6").Value = '625.45'
$ws.Range("E6").Value = '  -0.89%  '

# Row 7
$ws.Range("D7").Value = '0.393'
$ws.Range("E7").Value = '  +0.53%  '

# Row 8
$ws.Range("D8").Value = '0.702'
$ws.Range("E8").Value = '  +1.42%  '

# Row 9
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.19%  '

# Row 10
$ws.Range("D10").Value = '3.244.64'
$ws.Range("E10").Value = '  -1.14%  '

# Row 11
$ws.Range("D11").Value = '0.581'
$ws.Range("E11").Value = '  +0.04%  '

# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.179'
$ws.Range("E12").Value = '  -4.11%  '

# Row 13
$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").Value = '0.0000261'
$ws.Range("E13").Value = '  -0.93%  '

# Row 14
$ws.Range("B14").Value = 'WrappedBTC'
$ws.Range("C14").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D14").Value = '91.977.50'
$ws.Range("E14").Value = '  +3.89%  '

# Row 15
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '3.862.92'
$ws.Range("E15").Value = '  -0.63%  '

# Row 16
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").Value = '33.83'
$ws.Range("E16").Value = '  -0.97%  '

# Row 17
$ws.Range("D17").Value = '5.29'
$ws.Range("E17").Value = '  -1.70%  '

# Row 18
$ws.Range("D18").Value = '3.272.71'
$ws.Range("E18").Value = '  -0.49%  '

# Row 19
$ws.Range("D19").Value = '3.29'
$ws.Range("E19").Value = '  +5.66%  '

# Row 20
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '13.79'
$ws.Range("E20").Value = '  -2.37%  '

# Row 21
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = '441.29'
$ws.Range("E21").Value = '  +0.94%  '

# Row 22
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '8.77'
$ws.Range("E22").Value = '  -1.47%  '

# Row 23
$ws.Range("B23").Value = 'PEPE'
$ws.Range("C23").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D23").Value = '0.0000193'
$ws.Range("E23").Value = '  +43.69%  '

# Row 24
$ws.Range("D24").Value = '5.20'
$ws.Range("E24").Value = '  -3.68%  '

# Row 25
$ws.Range("D25").Value = '5.28'
$ws.Range("E25").Value = '  +1.37%  '

# Row 26
$ws.Range("D26").Value = '12.03'
$ws.Range("E26").Value = '  -2.57%  '

# Row 27
$ws.Range("D27").Value = '3.450.78'
$ws.Range("E27").Value = '  -0.03%  '

# Row 28
$ws.Range("D28").Value = '77.10'
$ws.Range("E28").Value = '  +0.03%  '

# Row 29
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.06%  '

# Row 30
$ws.Range("D30").Value = '0.171'
$ws.Range("E30").Value = '  -10.02%  '

# Row 31
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.15%  '

# Row 32
$ws.Range("B32").Value = 'dogwifhat'
$ws.Range("C32").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D32").Value = '4.19'
$ws.Range("E32").Value = '  +39.47%  '

# Row 33
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = '8.66'
$ws.Range("E33").Value = '  -2.74%  '

# Row 34
$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").Value = '546.95'
$ws.Range("E34").Value = '  -3.89%  '

# Row 35
$ws.Range("D35").Value = '7.00'
$ws.Range("E35").Value = '  -3.47%  '

# Row 36
$ws.Range("D36").Value = '1.92'
$ws.Range("E36").Value = '  -2.78%  '

# Row 37
$ws.Range("D37").Value = '1.27'
$ws.Range("E37").Value = '  -9.37%  '

# Row 38
$ws.Range("D38").Value = '22.42'
$ws.Range("E38").Value = '  -1.09%  '

# Row 39
$ws.Range("D39").Value = '22.46'
$ws.Range("E39").Value = '  +3.04%  '

# Row 40
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.13%  '

# Row 41
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '0.127'
$ws.Range("E41").Value = '  -8.02%  '

# Row 42
$ws.Range("D42").Value = '0.388'
$ws.Range("E42").Value = '  -2.85%  '

# Row 43
$ws.Range("D43").Value = '1.97'
$ws.Range("E43").Value = '  -2.72%  '

# Row 44
$ws.Range("E44").Value = '  -0.12%  '

# Row 45
$ws.Range("D45").Value = '149.63'
$ws.Range("E45").Value = '  -3.09%  '

# Row 46
$ws.Range("D46").Value = '45.07'
$ws.Range("E46").Value = '  -0.18%  '

# Row 47
$ws.Range("D47").Value = '177.22'
$ws.Range("E47").Value = '  -2.07%  '

# Row 48
$ws.Range("D48").Value = '0.126'
$ws.Range("E48").Value = '  +1.37%  '

# Row 49
$ws.Range("D49").Value = '1.27'
$ws.Range("E49").Value = '  -2.58%  '

# Row 50
$ws.Range("D50").Value = '0.630'
$ws.Range("E50").Value = '  +0.51%  '

# Row 51
$ws.Range("D51").Value = '4.17'
$ws.Range("E51").Value = '  -1.89%  '
